
$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($range) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Start -le $range.Start -and $range.End -le $p.Range.End) {
            return $i
        }
    }
    throw "No paragraph contains the given range"
}

function Find-Range($searchText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return $r
}

function Find-ParagraphIndex($searchText) {
    $r = Find-Range $searchText
    return Get-ParagraphIndexContaining $r
}

function Replace-Text($oldText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Replace-Text: text not found: $oldText"
    }
}

function Insert-ParagraphAfterIndex($idx, $text) {
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs($idx + 1)
    $newp.Range.InsertAfter($text)
    return $idx + 1
}

function Insert-ParagraphBeforeIndex($idx, $text) {
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphBefore()
    $newp = $d.Paragraphs($idx)
    $newp.Range.InsertAfter($text)
    return $idx
}

# --- Step 1: simple text replacements ---
Replace-Text ' Василь Пастушенко‑«Дмитро Яценко» (1941—1943).' ' Василь Пастушенко‑«Дмитро Яценко» (4.10.1941—1943).'
Replace-Text 'Іван Клим-«Митар».' 'Іван Клим-«Митар» (1941).'
Replace-Text ' Богдан Мовчан-«Степан Лисавка».' ' Богдан Мовчан-«Степан Лисавка» (4.10.1941—11.03.1943).'
Replace-Text ' Богдан Крицан -«Криця».' ' Богдан Крицан -«Криця» (4.10.1941—11.03.1943).'
Replace-Text 'Мережу підпілля ОУН(б) змогли створити члени підгрупи «Г» Південної похідної групи ОУН(б), які прибули до Запоріжжя після окупації німецькими військами 4 жовтня 1941 року. Перший склад обласного проводу мав такий вигляд: провідник – Василь Пастушенко («Дмитро Яценко»), організаційний референт – Іван Клим («Митар»), куратор Мелітопольського окружного проводу – Богдан Мовчан («Степан Лисавка»), інспектор діяльності низових ланок підпілля (та, імовірно, референт СБ) – Богдан Крицан («Криця»), Ілля (прізвище невідоме) – референт освіти і Степан Держко – член проводу (посада невідома).' 'Мережу підпілля ОУН(б) змогли створити члени підгрупи «Г» Південної похідної групи ОУН(б), які прибули до Запоріжжя після окупації німецькими військами 4 жовтня 1941 року. Перший склад обласного проводу мав такий вигляд: провідник — Василь Пастушенко («Дмитро Яценко»), організаційний референт — Іван Клим («Митар»), куратор Мелітопольського окружного проводу — Богдан Мовчан («Степан Лисавка»), інспектор діяльності низових ланок підпілля (та, імовірно, референт СБ) — Богдан Крицан («Криця»), Ілля (прізвище невідоме) — референт освіти і Степан Держко — член проводу (посада невідома).'

# --- Step 2: insert zapor1 paragraph after the "Merezhu pidpillia" paragraph ---
$idxMerezhu = Find-ParagraphIndex 'ще невідоме) — референт освіти і Степан Держко — член проводу (посада невідома).'
Insert-ParagraphAfterIndex $idxMerezhu 'zapor1' | Out-Null

# --- Step 3: insert zapor4 paragraph after "Na pochatok 1942 roku..." paragraph (before "Osnovni napriamky...") ---
$idx1942 = Find-ParagraphIndex 'Ще три райони підпорядковувалися обласному проводу напряму.'
Insert-ParagraphAfterIndex $idx1942 'zapor4' | Out-Null

# --- Step 4: split the "Pik represii..." paragraph into two, inserting zapor3 before and zapor2 between ---
$rSplit = Find-Range 'Від куль бандерівців загинуло троє співробітників СД. '
$idxPart1 = Get-ParagraphIndexContaining $rSplit
$rSplit.Collapse(0)
$rSplit.InsertParagraphAfter()
# zapor2 goes right after part1 (between part1 and part2)
Insert-ParagraphAfterIndex $idxPart1 'zapor2' | Out-Null
# zapor3 goes right before part1
Insert-ParagraphBeforeIndex $idxPart1 'zapor3' | Out-Null

# --- Step 5: append illustrations block after the second split part paragraph ---
$idxPart2 = Find-ParagraphIndex 'Того ж  місяця до в’язниці Мелітополя потрапив провідник Оси'
$idxCur = $idxPart2
$idxCur = Insert-ParagraphAfterIndex $idxCur ""
$idxCur = Insert-ParagraphAfterIndex $idxCur 'Ілюстрації:'
$idxCur = Insert-ParagraphAfterIndex $idxCur ""
$idxCur = Insert-ParagraphAfterIndex $idxCur '1.Провідник Запорізької області Василь Пастушенко - "Василь Ясенко" з мамою'
$idxCur = Insert-ParagraphAfterIndex $idxCur '2.Михайло Вінтонів - ”Михась” – провідник ОУН Мелітополmської округи Запорізької обл.'
$idxCur = Insert-ParagraphAfterIndex $idxCur '3.Костянтин Приходько - член ОУН(б) м. Бердянська Запорізької обл.'
$idxCur = Insert-ParagraphAfterIndex $idxCur '4.Схема виявлених, ліквідованих і розроблюваних НКДБ осередків ОУН у Запорізькій обл. на 20 червня 1945 року.'

Write-Host "Done. Total paragraphs:" $d.Paragraphs.Count